$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F3").Value = "SV Dessau 05"
$ws.Range("G3").Value = "1 FC BitterfeldWolfen"
$ws.Range("B3").Value = 6776470
$ws.Range("K3").Value = 2.2
